# Applies the "Substituindo login e cadastro por um modal" backlog update:
#  - status/progress updates on several tasks
#  - dependency values filled in for the "banco de dados" / "Integração" rows
#  - progress bump on "Criação da dropbox"
#  - a brand-new backlog row (31) for "Modelagem do banco de Dados"
#  - minor column-width + selection/view tweaks

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Progress (% Concluido) and Status updates -----------------------------
# Row 2: Criação da pag inicial -> 100% / Concluído
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = "Concluído"

# Row 3: Criação da tela de login -> 50% / Em progresso
$ws.Range("H3").Value = 0.5
$ws.Range("I3").Value = "Em progresso"

# Row 4: Criação da tela de cadastro -> 50% / Em progresso
$ws.Range("H4").Value = 0.5
$ws.Range("I4").Value = "Em progresso"

# Row 5: Criação da barra de pesquisa -> 25% / Em progresso
$ws.Range("H5").Value = 0.25
$ws.Range("I5").Value = "Em progresso"

# Row 9: Criação do footer -> 100% / Concluído
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = "Concluído"

# Row 10: Criação do botao -> 100% / Concluído
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = "Concluído"

# Row 12: Criação do formulário -> 50% / Em progresso
$ws.Range("H12").Value = 0.5
$ws.Range("I12").Value = "Em progresso"

# Row 14: Criação do card aulas -> 25% / Em progresso
$ws.Range("H14").Value = 0.25
$ws.Range("I14").Value = "Em progresso"

# Row 30: Criação da dropbox -> 50% (status stays "Em progresso")
$ws.Range("H30").Value = 0.5

# --- Dependency ("Depende do item") fills -----------------------------------
$ws.Range("C20").Value = 30
$ws.Range("C21").Value = 21

# --- New backlog row 31: Modelagem do banco de Dados ------------------------
$ws.Range("A31").Value = 30
$ws.Range("B31").Value = "Modelagem do banco de Dados"
$ws.Range("E31").Value = "Modelagem"
$ws.Range("F31").Value = "Dados"
$ws.Range("G31").Value = 8
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = "Não iniciado"

# --- Column width tweaks -----------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 19.285714285714285
$ws.Columns.Item(4).ColumnWidth = 23.571428571428573
$ws.Columns.Item(8).ColumnWidth = 14.571428571428571

# --- View / selection update --------------------------------------------------
$ws.Range("I16").Select()
